# Update the "Förändrad" (Changed) date column (C) for rows 2-12:
# increment the date serial value from 45174 to 45175 (one day later).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
